$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows ahead of the old "Incomplete:" row (currently row 33).
# This pushes the old row 33 ("Incomplete:") down to row 37, and the
# following rows (old 34-37) down to 38-41.
$ws.Rows("33:36").Insert()

# Row 33 stays fully blank (no cell data at all, like the other section
# spacer rows elsewhere in the sheet) - clear whatever formatting/content
# bled into it from the insert.
$ws.Rows("33").Clear()

# New "TV Guide" sub-heading (bold section header, same style family as the
# other bold headers such as "Incomplete:" / "1.0.2").
$ws.Range("A34").Value = "TV Guide"

# New "Done" item describing the skin setting for genre coloring. (Uses the
# same wrap/top-aligned, non-bold body style as the rest of the table - the
# row-insert above inherited the bold "heading" style, so un-bold it here.)
$ws.Range("A35").Value = "Provide setting to control genre coloring"
$ws.Range("A35").Font.Bold = $false
$ws.Range("A35").WrapText = $true
$ws.Range("A35").VerticalAlignment = -4160
$ws.Range("B35").Value = "Done"

# Blank spacer row (keeps a styled-but-empty A cell, matching row 12/22's
# pattern) between the new item and the "Incomplete:" heading.
$ws.Range("A36").Value = ""

# Move the view/selection the way the saved workbook shows it.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("A35").Select()
